$d = $word.ActiveDocument

# 1) Critério run: split the sentence from the formula with three line breaks
$d.Content.Find.Execute(
    "5,0 pontos.(Nota final",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5,0 pontos.^l^l^l(Nota final",
    2
)

# 2) Norma de recuperação run: split the sentence from the formula with two line breaks
$d.Content.Find.Execute(
    "sua nota final.(Nota final",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sua nota final.^l^l(Nota final",
    2
)

# 3) Bibliografia run: split the two references with two line breaks
$d.Content.Find.Execute(
    "1473p.B)OGA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1473p.^l^lB)OGA",
    2
)
